# Fruta / hortaliza, semanal
# Insert two new weekly price records (rows) right before the current row 16,
# pushing the existing data down by two rows, then populate the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 16 (existing rows 16.. shift down to 18..)
$ws.Rows("16:17").Insert()

function Set-RowValues($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Unidad, $Origen, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value = 1
    $ws.Cells.Item($Row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($Row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 15
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100103
    $ws.Cells.Item($Row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($Row, 9).Value = 100103006
    $ws.Cells.Item($Row, 10).Value = "Nectarín"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 15).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 16).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = 18
}

# New row 16: Artic Snow, Segunda, 2022-02-28
Set-RowValues 16 44620 "Artic Snow" "Segunda" 270 19000 20000 19500 "`$/bandeja 18 kilos granel" "Región de O'Higgins" 1083

# New row 17: August Red, Segunda, 2022-02-28
Set-RowValues 17 44620 "August Red" "Segunda" 250 19000 20000 19500 "`$/bandeja 18 kilos granel" "Región de O'Higgins" 1083

Write-Output "edit complete"
